# NIT-9016428030.xlsx - "Estado de Cuenta" update
#
# The "Periodo Mora" detail table (rows 16-21, columns B:J on sheet "Hoja1")
# needs its data rows re-ordered so the periods run in ascending order
# (2412, 2502, 2503, 2504, 2505, 2506) instead of the previous
# (2505, 2504, 2503, 2502, 2412, 2506) order. Sorting the block by the
# "Periodo Mora" column (E) carries every other column (Valor Mora in F,
# Salario Basico in G, etc.) along with each row, which is exactly how the
# two swapped "Valor Mora" figures (56940 / 52000) end up back on the right
# period row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$dataRange = $ws.Range("B16:J21")
$sortKey = $ws.Range("E16:E21")

$dataRange.Sort($sortKey, 1)
